$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the last existing row (99) with revised quarterly figures
$ws.Range("B99").Value = 23.8
$ws.Range("C99").Value = 19.2
$ws.Range("D99").Value = 2.4
$ws.Range("E99").Value = 21.5
$ws.Range("F99").Value = -0.3
$ws.Range("G99").Value = 25.5
$ws.Range("H99").Value = 20.5

# Append a new row (100) for the next quarterly period "01-04-2021"
# Force the cell to remain text so Excel does not auto-convert the
# date-like string into a serial date number, then restore the default
# (General) formatting so no residual style is left behind.
$ws.Range("A100").NumberFormat = "@"
$ws.Range("A100").Value = "01-04-2021"
$ws.Range("A100").Style = "Normal"

$ws.Range("B100").Value = 21.8
$ws.Range("C100").Value = 20.2
$ws.Range("D100").Value = 3.3
$ws.Range("E100").Value = 18.5
$ws.Range("F100").Value = 0.2
$ws.Range("G100").Value = 22.4
$ws.Range("H100").Value = 20.7
$ws.Range("I100").Value = 0.2
